$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-point the "last row" special border styling.
#    Row 535 is currently the last data row and carries the special
#    bottom-border style-set. Row 549 will become the new last row, so move
#    that style-set over to row 549 (A:M) first, then restore row 535 (A:M)
#    to the regular alternating (odd-row) style used throughout the table.
# ---------------------------------------------------------------------------
$ws.Range("A535:M535").Copy()
$ws.Range("A549:M549").PasteSpecial(-4122)

$ws.Range("A533:M533").Copy()
$ws.Range("A535:M535").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Stamp the normal alternating row styles (format only) onto the newly
#    added rows 536-548, picking the template row that matches both the
#    target row's odd/even parity and whether it needs an M-column value
#    (first sub-answer) or an N-column value (second sub-answer).
#      - even row, N column -> template row 4
#      - odd  row, N column -> template row 3
#      - even row, M column -> template row 6
#      - odd  row, M column -> template row 5
# ---------------------------------------------------------------------------
$nEvenTemplate = "A4:N4"
$nOddTemplate  = "A3:N3"
$mEvenTemplate = "A6:N6"
$mOddTemplate  = "A5:N5"

$ws.Range($nEvenTemplate).Copy()
$ws.Range("A536:N536").PasteSpecial(-4122)

$ws.Range($nOddTemplate).Copy()
$ws.Range("A537:N537").PasteSpecial(-4122)

$ws.Range($nEvenTemplate).Copy()
$ws.Range("A538:N538").PasteSpecial(-4122)

$ws.Range($mOddTemplate).Copy()
$ws.Range("A539:N539").PasteSpecial(-4122)

$ws.Range($mEvenTemplate).Copy()
$ws.Range("A540:N540").PasteSpecial(-4122)

$ws.Range($mOddTemplate).Copy()
$ws.Range("A541:N541").PasteSpecial(-4122)

$ws.Range($nEvenTemplate).Copy()
$ws.Range("A542:N542").PasteSpecial(-4122)

$ws.Range($mOddTemplate).Copy()
$ws.Range("A543:N543").PasteSpecial(-4122)

$ws.Range($nEvenTemplate).Copy()
$ws.Range("A544:N544").PasteSpecial(-4122)

$ws.Range($mOddTemplate).Copy()
$ws.Range("A545:N545").PasteSpecial(-4122)

$ws.Range($mEvenTemplate).Copy()
$ws.Range("A546:N546").PasteSpecial(-4122)

$ws.Range($nOddTemplate).Copy()
$ws.Range("A547:N547").PasteSpecial(-4122)

$ws.Range($mEvenTemplate).Copy()
$ws.Range("A548:N548").PasteSpecial(-4122)

# Row 549 (new last row) already carries the old row-535 last-row style-set
# on A:M. Give it the matching N-column last-row treatment by starting from
# the regular N-column style (row 531) and darkening the bottom border to
# the same "last row" accent color used by the other last-row columns.
$ws.Range("N531").Copy()
$ws.Range("N549").PasteSpecial(-4122)
$ws.Range("N549").Borders.Item(9).Color = 6631236

# ---------------------------------------------------------------------------
# 3) Write the 14 new response rows (536-549).
# ---------------------------------------------------------------------------
$ws.Range("A536").Value = 45570.380015821764
$ws.Range("B536").Value = "min010417@gmail.com"
$ws.Range("C536").Value = "환경생명공학과"
$ws.Range("D536").Value = 20203702
$ws.Range("E536").Value = "강채민"
$ws.Range("F536").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G536").Value = 0.1
$ws.Range("H536").Value = "6:4"
$ws.Range("I536").Value = "20분의 1"
$ws.Range("J536").Value = "20만호, 69만명"
$ws.Range("K536").Value = "충청"
$ws.Range("L536").Value = "Black"
$ws.Range("N536").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A537").Value = 45570.399496157406
$ws.Range("B537").Value = "tlsdmsco1130@naver.com"
$ws.Range("C537").Value = "미디어스쿨"
$ws.Range("D537").Value = 20242532
$ws.Range("E537").Value = "신은채"
$ws.Range("F537").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G537").Value = 0.9
$ws.Range("H537").Value = "6:4"
$ws.Range("I537").Value = "20분의 1"
$ws.Range("J537").Value = "20만호, 69만명"
$ws.Range("K537").Value = "평안"
$ws.Range("L537").Value = "Black"
$ws.Range("N537").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A538").Value = 45570.492105428246
$ws.Range("B538").Value = "kyj57980@gmail.com"
$ws.Range("C538").Value = "사회복지학과"
$ws.Range("D538").Value = 20202319
$ws.Range("E538").Value = "김예진"
$ws.Range("F538").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G538").Value = 0.1
$ws.Range("H538").Value = "3:7"
$ws.Range("I538").Value = "10분의 1"
$ws.Range("J538").Value = "44만호, 153만명"
$ws.Range("K538").Value = "전라"
$ws.Range("L538").Value = "Black"
$ws.Range("N538").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A539").Value = 45570.58668759259
$ws.Range("B539").Value = "a01075976680@gmail.com"
$ws.Range("C539").Value = "반도체디스플레이스쿨"
$ws.Range("D539").Value = 20243309
$ws.Range("E539").Value = "김우진"
$ws.Range("F539").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G539").Value = 0.3
$ws.Range("H539").Value = "3:7"
$ws.Range("I539").Value = "10분의 1"
$ws.Range("J539").Value = "44만호, 153만명"
$ws.Range("K539").Value = "전라"
$ws.Range("L539").Value = "Red"
$ws.Range("M539").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A540").Value = 45570.60545430555
$ws.Range("B540").Value = "sohn1118@naver.com"
$ws.Range("C540").Value = "생명과학과"
$ws.Range("D540").Value = 20223519
$ws.Range("E540").Value = "손정빈"
$ws.Range("F540").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G540").Value = 0.1
$ws.Range("H540").Value = "6:4"
$ws.Range("I540").Value = "20분의 1"
$ws.Range("J540").Value = "20만호, 69만명"
$ws.Range("K540").Value = "충청"
$ws.Range("L540").Value = "Red"
$ws.Range("M540").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A541").Value = 45570.668387881946
$ws.Range("B541").Value = "hlu20242513@gmail.com"
$ws.Range("C541").Value = "미디어스쿨"
$ws.Range("D541").Value = 20242513
$ws.Range("E541").Value = "김예준"
$ws.Range("F541").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G541").Value = 0.3
$ws.Range("H541").Value = "6:4"
$ws.Range("I541").Value = "10분의 1"
$ws.Range("J541").Value = "15만호,  32만명"
$ws.Range("K541").Value = "경상"
$ws.Range("L541").Value = "Red"
$ws.Range("M541").Value = "반대한다."

$ws.Range("A542").Value = 45570.68802954861
$ws.Range("B542").Value = "hhy062700@naver.com"
$ws.Range("C542").Value = "언어청각학부"
$ws.Range("D542").Value = 20243971
$ws.Range("E542").Value = "황희영"
$ws.Range("F542").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G542").Value = 0.9
$ws.Range("H542").Value = "6:4"
$ws.Range("I542").Value = "10분의 1"
$ws.Range("J542").Value = "44만호, 153만명"
$ws.Range("K542").Value = "평안"
$ws.Range("L542").Value = "Black"
$ws.Range("N542").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A543").Value = 45570.6892605787
$ws.Range("B543").Value = "lapter1122@naver.com"
$ws.Range("C543").Value = "광고홍보학과"
$ws.Range("D543").Value = 20232635
$ws.Range("E543").Value = "진유진"
$ws.Range("F543").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G543").Value = 0.1
$ws.Range("H543").Value = "6:4"
$ws.Range("I543").Value = "20분의 1"
$ws.Range("J543").Value = "20만호, 69만명"
$ws.Range("K543").Value = "충청"
$ws.Range("L543").Value = "Red"
$ws.Range("M543").Value = "반대한다."

$ws.Range("A544").Value = 45570.701951423616
$ws.Range("B544").Value = "psh020509@naver.com"
$ws.Range("C544").Value = "경제학과 "
$ws.Range("D544").Value = 20212820
$ws.Range("E544").Value = "박시환"
$ws.Range("F544").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G544").Value = 0.1
$ws.Range("H544").Value = "6:4"
$ws.Range("I544").Value = "20분의 1"
$ws.Range("J544").Value = "20만호, 69만명"
$ws.Range("K544").Value = "충청"
$ws.Range("L544").Value = "Black"
$ws.Range("N544").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A545").Value = 45570.71003824074
$ws.Range("B545").Value = "rhoy3156@hallym.ac.kr"
$ws.Range("C545").Value = "법학과"
$ws.Range("D545").Value = 20202719
$ws.Range("E545").Value = "노원철"
$ws.Range("F545").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G545").Value = 0.1
$ws.Range("H545").Value = "6:4"
$ws.Range("I545").Value = "10분의 1"
$ws.Range("J545").Value = "44만호, 153만명"
$ws.Range("K545").Value = "충청"
$ws.Range("L545").Value = "Red"
$ws.Range("M545").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A546").Value = 45570.72153730324
$ws.Range("B546").Value = "youu0729@naver.com"
$ws.Range("C546").Value = "사회복지학부"
$ws.Range("D546").Value = 20242335
$ws.Range("E546").Value = "유이현"
$ws.Range("F546").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G546").Value = 0.1
$ws.Range("H546").Value = "6:4"
$ws.Range("I546").Value = "20분의 1"
$ws.Range("J546").Value = "20만호, 69만명"
$ws.Range("K546").Value = "충청"
$ws.Range("L546").Value = "Red"
$ws.Range("M546").Value = "모름/무응답"

$ws.Range("A547").Value = 45570.73606974537
$ws.Range("B547").Value = "jwtp724@naver.com"
$ws.Range("C547").Value = "콘텐츠IT"
$ws.Range("D547").Value = 20206504
$ws.Range("E547").Value = "박이선"
$ws.Range("F547").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G547").Value = 0.3
$ws.Range("H547").Value = "5:5"
$ws.Range("I547").Value = "20분의 1"
$ws.Range("J547").Value = "15만호,  32만명"
$ws.Range("K547").Value = "전라"
$ws.Range("L547").Value = "Black"
$ws.Range("N547").Value = "찬성한다."

$ws.Range("A548").Value = 45570.73702164352
$ws.Range("B548").Value = "dldpwls5245@naver.com"
$ws.Range("C548").Value = "법학과"
$ws.Range("D548").Value = 20182747
$ws.Range("E548").Value = "이예진"
$ws.Range("F548").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G548").Value = 0.1
$ws.Range("H548").Value = "4:6"
$ws.Range("I548").Value = "10분의 1"
$ws.Range("J548").Value = "44만호, 153만명"
$ws.Range("K548").Value = "평안"
$ws.Range("L548").Value = "Red"
$ws.Range("M548").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A549").Value = 45570.737177800926
$ws.Range("B549").Value = "sehyun1901@gmail.com"
$ws.Range("C549").Value = "경영대학"
$ws.Range("D549").Value = 20243020
$ws.Range("E549").Value = "임세현"
$ws.Range("F549").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G549").Value = 0.3
$ws.Range("H549").Value = "6:4"
$ws.Range("I549").Value = "15분의 1"
$ws.Range("J549").Value = "44만호, 153만명"
$ws.Range("K549").Value = "전라"
$ws.Range("L549").Value = "Black"
$ws.Range("N549").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# ---------------------------------------------------------------------------
# 4) Grow the "Form_Responses1" table so it covers the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N549"))
